$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear out the old "on-board parts" rows 8-17 (columns A-F) so stale
#     cells/styles don't linger, then rebuild them with the new BOM content. ---
$ws.Range("A8:F17").Clear()

# Row 8: was "1k resistor"/"541-3991-1-ND" -> now "0.1uF ceramic", qty 3
$ws.Range("A8").Value = "0.1uF ceramic"
$ws.Range("C8").Value = 3

# Row 9: was "10k resistor" -> now "220 ohm resistor", qty 2
$ws.Range("A9").Value = "220 ohm resistor"
$ws.Range("C9").Value = 2

# Row 10: was "reset button"/"401-1426-1-ND" -> now "330 ohm resistor", qty 5
$ws.Range("A10").Value = "330 ohm resistor"
$ws.Range("C10").Value = 5

# Row 11: new "10k resistor" row (moved down from the old row 9 slot)
$ws.Range("A11").Value = "10k resistor"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 0
$ws.Range("F11").Value = "Have"

# Row 12: "reset button" row (moved down from the old row 10 slot)
$ws.Range("A12").Value = "reset button"
$ws.Range("B12").Value = "401-1426-1-ND"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = 0.52

# Row 13: "Level shifter" (moved down from old row 12)
$ws.Range("A13").Value = "Level shifter"
$ws.Range("B13").Value = "296-12163-1-ND"
$ws.Range("C13").Value = 1
$ws.Range("E13").Value = 0.43

# Row 14: was "SMA connector" -> now "SMA PCB connector" (moved down from old row 13)
$ws.Range("A14").Value = "SMA PCB connector"
$ws.Range("B14").Value = "A97594-ND"
$ws.Range("C14").Value = 1
$ws.Range("E14").Value = 2.17

# Row 15-17: brand-new rows describing the new SMA / resistor parts
$ws.Range("A15").Value = "SMA right-angle panel mount connector"
$ws.Range("A16").Value = "SMA wire connectors"
$ws.Range("A17").Value = "50 ohm SMA coax"
$ws.Range("D17").Value = 0
$ws.Range("F17").Value = "Have"

# Row 15 wraps onto two lines like row 10 already does, so give it the same height.
$ws.Rows.Item(15).RowHeight = 28.8

# B10 keeps the small Arial "note" formatting used elsewhere (e.g. B7) even
# though it no longer holds a value.
$ws.Range("B7").Copy()
$ws.Range("B10").PasteSpecial(-4122)

# Restore the current selection to E7.
$ws.Range("E7").Select()
